$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 1..18) holds numeric-looking text values. Per the target
# diff only rows 3-10, 17 and 18 actually change their displayed value;
# rows 1, 2 and 11-16 must stay exactly as they are.
$changes = @(
    , @(3,  "47.3")
    , @(4,  "56.4")
    , @(5,  "26.6")
    , @(6,  "69.9")
    , @(7,  "43.4")
    , @(8,  "62.4")
    , @(9,  "50.2")
    , @(10, "12.9")
    , @(17, "166.2")
    , @(18, "66.4")
)

foreach ($pair in $changes) {
    $row = $pair[0]
    $text = $pair[1]
    $cell = $ws.Cells.Item($row, 1)

    # These values look numeric, so a plain `.Value = $text` assignment
    # would be auto-coerced into a real number by Excel. The source
    # workbook stores them as shared-string text, so force text storage
    # via NumberFormat "@" and then clear the format change that
    # introduces, restoring the cell to its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}
